$d = $word.ActiveDocument

function Replace-DocText($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "WARNING: text not found:" $findText
    }
}

# --- Title & author name -------------------------------------------------
Replace-DocText "Unveiling the Enigma of Dark Matter" "Unveiling the Enigmatic Symphony of Biology"
Replace-DocText " Isabella Rosselli" " Sarah Williams"

# --- Author email (paragraph 3): split into dr + . + sarah25@academics ---
$emailPara = $d.Paragraphs.Item(3)
$emailParaText = $emailPara.Range.Text
$oldUser = "irosselli@physics"
$userOffset = $emailParaText.IndexOf($oldUser)
$userStart = $emailPara.Range.Start + $userOffset
$userEnd = $userStart + $oldUser.Length
$userRange = $d.Range($userStart, $userEnd)
$userRange.Text = "dr"
# insert a separate '.' run, then the new username text, right after 'dr'
$afterDr = $d.Range($userStart + 2, $userStart + 2)
$afterDr.InsertAfter(".")
$dotEnd = $userStart + 2 + 1
$afterDot = $d.Range($dotEnd, $dotEnd)
$afterDot.InsertAfter("sarah25@academics")

# --- Main body paragraph (paragraph 5): replace each of the three ---
# --- double-line-break-separated segments in turn                  ---
Replace-DocText "Within the vast cosmic tapestry, there lies an elusive substance known as dark matter, a mysterious entity that permeates the universe, exerting a gravitational influence yet remaining invisible to our most powerful telescopes. Its existence is inferred through its gravitational effects on visible matter, such as stars and galaxies, and its enigmatic nature has captivated the imaginations of scientists and laypeople alike. In this exploration, we will delve into the enigma of dark matter, examining its properties, potential candidates, and the ongoing quest to unveil its true identity." "In the vast auditorium of life, Biology stands as a mesmerizing drama, unfolding a grand narrative of existence. It's a symphony of intricate processes, a cosmic dance of molecules and cells that weaves the tapestry of life. Behold the microscopic realm, a universe teeming with activity, where organisms engage in a delicate ballet of survival. Behold the grandeur of evolution, the chronicle of life's relentless journey through epochs, adapting and diversifying in response to the relentless rhythm of change, leaving an indelible mark on Earth's geological memoirs. Biology, a kaleidoscope of colors and shapes, offers a breathtaking glimpse into the enigmatic symphony of life, inviting us to unravel its secrets. Unveil the intricate mechanisms that govern heredity, the blueprint of life encoded within the DNA molecule, a script that orchestrates the construction of every living being."
Replace-DocText "Dark matter's presence is evident in the intricate dance of galaxies, where its gravitational pull shapes their motions and influences their structures. Observations have revealed that the mass of galaxies, as inferred from their gravitational effects, often far exceeds the mass of the visible matter they contain. This discrepancy points to the existence of a hidden mass component, an unseen force that governs the dynamics of the cosmos. Additionally, gravitational lensing, the bending of light around massive objects, provides further evidence for the existence of dark matter, as the observed distortions in the shapes of distant galaxies align with the predictions of its presence." "Journey into the depths of cellular biology, where organelles perform a synchronized symphony, executing vital functions that sustain life. Explore the intricate dance of photosynthesis, where plants capture the sun's radiant energy, orchestrating a symphony of biochemical transformations that nourish the planet. Delve into the intricacies of the human body, a marvel of engineering, where systems harmoniously collaborate, a symphony of physiological processes that maintain our existence. Biology, the study of life, is a boundless realm of discovery, an invitation to explore the very essence of existence. It beckons us to decode the enigmatic symphony of life, a tapestry woven by the hand of evolution, holding the key to our understanding of the universe and our place within it."
Replace-DocText "The nature of dark matter remains shrouded in mystery, with various theories attempting to unravel its composition. One leading candidate is Weakly Interacting Massive Particles (WIMPs), hypothetical particles that are massive but interact with ordinary matter only through weak nuclear forces, making them difficult to detect directly. Another possibility is that dark matter consists of primordial black holes, formed in the early universe through the collapse of massive clouds of gas. However, these candidates have faced challenges in explaining all the observed properties of dark matter, leaving the true nature of this enigmatic substance still unresolved." "Comprehending Biology is embarking on an odyssey of exploration, a quest to understand the symphony of life in all its splendor. It's about unraveling the secrets of the natural world. Understanding the delicate balance of ecosystems, the interplay of organisms, and the intricate mechanisms that govern the harmony of life. Biology unveils the wonders of biodiversity, showcasing the astonishing array of species that inhabit our planet. Through its study, we grasp the importance of preserving the fragile equilibrium of the environment. Biology invites us to ponder the profound questions of existence, the origin of life, the nature of consciousness, and the intricate relationship between the living and nonliving world. Beyond its intellectual pursuits, Biology cultivates an appreciation for the interconnectedness of all life, fostering a sense of awe and wonder at the grandeur of the natural world. It empowers us with knowledge that can shape our decisions and actions, creating a more sustainable and harmonious co-existence with our environment."

# --- Move lastRenderedPageBreak from Summary paragraph into the body ---
# --- paragraph, right before 'and the intricate relationship...'     ---
$markerText = "the nature of consciousness, and the intricate relationship"
$bodyRange = $d.Content
$found = $bodyRange.Find.Execute($markerText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitText = "the nature of consciousness, "
    $insertPos = $bodyRange.Start + $splitText.Length
    $insertRange = $d.Range($insertPos, $insertPos)
    $insertRange.InsertAfter([char]2)
}

# --- Summary heading paragraph body text ---
Replace-DocText "Dark matter stands as one of the most perplexing enigmas in modern physics. Its gravitational influence shapes the universe, yet its true nature eludes our understanding. The search for dark matter particles continues, with experiments and observations probing the depths of the cosmos in pursuit of clues to its identity. Unveiling the secrets of dark matter holds the promise of revolutionizing our understanding of the universe, providing insights into the fundamental forces that govern its vast expanse." "Biology unravels the captivating tapestry of life, revealing the symphony of processes that govern the existence of organisms, from the smallest molecules to the grandest ecosystems. Its exploration unveils the intricacies of heredity, the mechanisms of cellular function, and the wonders of diversity among species, showcasing the intricate balance of life on Earth. Biology nurtures an appreciation for the interconnectedness of all living beings and inspires us to preserve the delicate harmony of our planet. Through its study, we gain an understanding of the natural world, shaping our decisions and actions towards a more sustainable and harmonious co-existence with our environment."

# --- Append trailing empty paragraph at the very end of the document ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

